$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ECs', 'L1cam', 'Itgav', 'ECs', 3, 1, 19.72083766666667, 59.162513, 0.8016210077351786, 0.8016210077351787, 3, 1, 16.535604, 49.606812, 0.2120453146491552, 0.2120453146491552, 326.095962204284, 2934.863659838556, 0.1699799788145789, 0.1699799788145789),
    @('ECs', 'L1cam', 'Itgav', 'sCs', 3, 1, 19.72083766666667, 59.162513, 0.8016210077351786, 0.8016210077351787, 3, 1, 40.62063066666667, 121.861892, 0.5209011059384622, 0.5209011059384622, 801.0728632949553, 7209.655769654597, 0.4175652694727591, 0.4175652694727591),
    @('ECs', 'L1cam', 'Itgav', 'FAPs', 3, 1, 19.72083766666667, 59.162513, 0.8016210077351786, 0.8016210077351787, 3, 1, 20.825229, 62.475687, 0.2670535794123827, 0.2670535794123827, 410.690960480159, 3696.218644321431, 0.2140757594478407, 0.2140757594478408),
    @('sCs', 'L1cam', 'Itgav', 'ECs', 1, 0.3333333333333333, 0.099159, 0.297477, 0.004030657259573097, 0.004030657259573097, 3, 1, 16.535604, 49.606812, 0.2120453146491552, 0.2120453146491552, 1.639653957036, 14.756885613324, 0.0008546819868490789, 0.0008546819868490793),
    @('sCs', 'L1cam', 'Itgav', 'sCs', 1, 0.3333333333333333, 0.099159, 0.297477, 0.004030657259573097, 0.004030657259573097, 3, 1, 40.62063066666667, 121.861892, 0.5209011059384622, 0.5209011059384622, 4.027901116276, 36.251110046484, 0.002099573824170517, 0.002099573824170518),
    @('sCs', 'L1cam', 'Itgav', 'FAPs', 1, 0.3333333333333333, 0.099159, 0.297477, 0.004030657259573097, 0.004030657259573097, 3, 1, 20.825229, 62.475687, 0.2670535794123827, 0.2670535794123827, 2.065008882411, 18.585079941699, 0.001076401448553501, 0.001076401448553501),
    @('FAPs', 'L1cam', 'Itgav', 'ECs', 3, 1, 4.781202, 14.343606, 0.1943483350052483, 0.1943483350052483, 3, 1, 16.535604, 49.606812, 0.2120453146491552, 0.2120453146491552, 79.06006291600801, 711.540566244072, 0.04121065384772729, 0.0412106538477273),
    @('FAPs', 'L1cam', 'Itgav', 'sCs', 3, 1, 4.781202, 14.343606, 0.1943483350052483, 0.1943483350052483, 3, 1, 40.62063066666667, 121.861892, 0.5209011059384622, 0.5209011059384622, 194.215440584728, 1747.938965262552, 0.1012362626415326, 0.1012362626415326),
    @('FAPs', 'L1cam', 'Itgav', 'FAPs', 3, 1, 4.781202, 14.343606, 0.1943483350052483, 0.1943483350052483, 3, 1, 20.825229, 62.475687, 0.2670535794123827, 0.2670535794123827, 99.569626545258, 896.1266389073221, 0.05190141851598842, 0.05190141851598842)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}